# Auto - Update data with bot!
# Updates title/link pairs for a few rows in the blog-used list sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Minho lee / lumiamitie blog entry
$ws.Range("D3").Value = "R로 프로덕션 개발 운영 환경 구축하기"
$ws.Range("E3").Value = "https://lumiamitie.github.io/data/r-for-production/"

# Row 5: Gongdols math / angeloyeo blog entry
$ws.Range("D5").Value = "양의 정부호 행렬"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2021/12/20/positive_definite.html"

# Row 9: 파비블로그 / pabii blog entry
$ws.Range("D9").Value = "MBA AI/BigData 2nd term 시험 문제 공개"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/mba-ai-bigdata-2nd-term-exam-sample-1/#utm_source=rss&utm_medium=rss&utm_campaign=mba-ai-bigdata-2nd-term-exam-sample-1"
